$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.657.87"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.677.06"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'314.12"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.3937"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").Value = "'0.3951"
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("D9").Value = "'1.000"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'1.408"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").Value = "'50.95"
$ws.Range("E11").Value = "  -5.21%  "
$ws.Range("D12").Value = "'0.08661"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "'25.30"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "'7.349"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "'0.00001320"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "'7.729"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "1.676.73"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "'94.04"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").Value = "'21.14"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "'7.094"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'13.96"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "24.667.47"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").Value = "'2.359"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").Value = "'2.793"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("D27").Value = "'23.08"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'5.879"
$ws.Range("E28").Value = "  -10.61%  "
$ws.Range("D29").Value = "'160.13"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'147.03"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").Value = "'8.354"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").Value = "'2.496"
$ws.Range("E32").Value = "  +10.14%  "
$ws.Range("D33").Value = "1.856.75"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "'0.03096"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'0.08315"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("D36").Value = "'6.972"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").Value = "'0.2810"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").Value = "'0.9947"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "'0.09545"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "'1.519"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").Value = "'10.33"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").Value = "'0.7933"
$ws.Range("E42").Value = "  -6.90%  "
$ws.Range("D43").Value = "'13.56"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "'16.65"
$ws.Range("E44").Value = "  -6.60%  "
$ws.Range("D45").Value = "'0.7126"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").Value = "'2.567"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("D47").Value = "'4.170"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "'0.08658"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'1.332"
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("D51").Value = "'137.83"
$ws.Range("E51").Value = "  -2.23%  "
